# Apply "contingencies with rene fine" edit:
# - Extend the table from columns A:O to A:Q (add columns P and Q)
# - Row 1 (header row): P1 = 14, Q1 = 15, formatted like the rest of the header row
# - Rows 2-25 (data rows):
#     * swap values in columns I and K (1<->2)
#     * swap values in columns M and O (1<->2)
#     * add new columns P = 2, Q = 2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): fill new P1/Q1 and copy formatting from O1 ---
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Data rows 2-25 ---
for ($r = 2; $r -le 25; $r++) {
    # Swap I <-> K
    $iVal = $ws.Cells.Item($r, 9).Value2   # column I
    $kVal = $ws.Cells.Item($r, 11).Value2  # column K
    $ws.Cells.Item($r, 9).Value = $kVal
    $ws.Cells.Item($r, 11).Value = $iVal

    # Swap M <-> O
    $mVal = $ws.Cells.Item($r, 13).Value2  # column M
    $oVal = $ws.Cells.Item($r, 15).Value2  # column O
    $ws.Cells.Item($r, 13).Value = $oVal
    $ws.Cells.Item($r, 15).Value = $mVal

    # New columns P and Q
    $ws.Cells.Item($r, 16).Value = 2   # column P
    $ws.Cells.Item($r, 17).Value = 2   # column Q
}
